$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 93..188) down by two rows (to 95..190) in
# a single Copy/Paste so Excel's own engine handles the row-by-row relocation
# (including carrying the date-format style on column D) without us clobbering
# source rows while still reading them.
$ws.Range("A93:R188").Copy($ws.Range("A95:R190"))

# Row 93 gets new data (Primera / Paine, $/kilo).
$ws.Cells.Item(93, 4).Value = 44895   # D93 Fecha
$ws.Cells.Item(93, 10).Value = 2000   # J93 Volumen
$ws.Cells.Item(93, 11).Value = 400    # K93 Precio minimo
$ws.Cells.Item(93, 12).Value = 400    # L93 Precio maximo
$ws.Cells.Item(93, 13).Value = 400    # M93 Precio promedio ponderado
$ws.Cells.Item(93, 15).Value = "Paine" # O93 Origen
$ws.Cells.Item(93, 16).Value = 400    # P93 Precio $/Kg

# Row 94 gets new data (Primera / Peru, $/kilo).
$ws.Cells.Item(94, 4).Value = 44895   # D94 Fecha
$ws.Cells.Item(94, 9).Value = "Primera" # I94 Calidad
$ws.Cells.Item(94, 10).Value = 2000   # J94 Volumen
$ws.Cells.Item(94, 11).Value = 500    # K94 Precio minimo
$ws.Cells.Item(94, 12).Value = 500    # L94 Precio maximo
$ws.Cells.Item(94, 13).Value = 500    # M94 Precio promedio ponderado
$ws.Cells.Item(94, 14).Value = "$/kilo" # N94 Unidad de comercializacion
$ws.Cells.Item(94, 15).Value = "Perú" # O94 Origen
$ws.Cells.Item(94, 16).Value = 500    # P94 Precio $/Kg
